$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.0013667345047
$ws.Range("B1").Value = 1.717131495475769
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.134654998779297
$ws.Range("E1").Value = 1.242484092712402
